# "Some enhancements for Page Load"
# Updates the CHROME_0218.. placeholder ERP/login identifiers to CHROME_0219..
# across the three affected sheets, and refreshes the selected
# cell/active-sheet bookkeeping left behind by the author's navigation.

$wb = $excel.ActiveWorkbook

# --- loginTest -----------------------------------------------------------
$login = $wb.Worksheets.Item("loginTest")
$login.Range("B2").Value = "CHROME_021916aa"

# --- accountCreatonTest ---------------------------------------------------
$acct = $wb.Worksheets.Item("accountCreatonTest")

$acct.Range("D2").Value  = "CHROME_021916aa"
$acct.Range("H2").Value  = "CHROME_021916aa"

$acct.Range("D3").Value  = "CHROME_021916ab"
$acct.Range("H3").Value  = "CHROME_021916ab"

$acct.Range("D4").Value  = "CHROME_021916ac"
$acct.Range("H4").Value  = "CHROME_021916ac"

$acct.Range("D5").Value  = "CHROME_021916ad"
$acct.Range("H5").Value  = "CHROME_021916ad"

$acct.Range("D6").Value  = "CHROME_021916ae"
$acct.Range("H6").Value  = "CHROME_021916ae"

$acct.Range("D7").Value  = "CHROME_021916af"
$acct.Range("H7").Value  = "CHROME_021916af"

$acct.Range("D8").Value  = "CHROME_021916ag"
$acct.Range("H8").Value  = "CHROME_021916ag"

$acct.Range("D9").Value  = "CHROME_021916ah"
$acct.Range("H9").Value  = "CHROME_021916ah"

$acct.Range("D10").Value = "CHROME_021916ai"
$acct.Range("H10").Value = "CHROME_021916ai"

$acct.Range("D11").Value = "CHROME_021916aj"
$acct.Range("H11").Value = "CHROME_021916aj"

# --- profileUpdateTest -----------------------------------------------------
$profile = $wb.Worksheets.Item("profileUpdateTest")
$profile.Range("E2").Value = "CHROME_021916aa"
$profile.Range("E3").Value = "CHROME_021916aa"

# --- Selection / active sheet bookkeeping ---------------------------------
# Replay the cursor moves the author made before saving: loginTest and
# accountCreatonTest selections moved, and profileUpdateTest became the
# active (visible-on-open) sheet.
$login.Activate()
$login.Range("B7").Select()

$acct.Activate()
$acct.Range("H11").Select()

$profile.Activate()
$profile.Range("B3").Select()
